$wb = $excel.ActiveWorkbook

# Citywide Totals (sheet 1)
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 9).Value = 2798
$ws.Cells.Item(3, 9).Value = 2906
$ws.Cells.Item(4, 5).Value = 1963
$ws.Cells.Item(4, 9).Value = 707
$ws.Cells.Item(5, 9).Value = 258
$ws.Cells.Item(6, 9).Value = 3317
$ws.Cells.Item(7, 5).Value = 25967
$ws.Cells.Item(7, 9).Value = 9986

# By Neighborhood (sheet 2)
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(7, 9).Value = 335
$ws.Cells.Item(8, 9).Value = 635
$ws.Cells.Item(14, 9).Value = 51
$ws.Cells.Item(19, 9).Value = 274
$ws.Cells.Item(20, 9).Value = 252
$ws.Cells.Item(23, 9).Value = 87
$ws.Cells.Item(27, 9).Value = 90
$ws.Cells.Item(29, 9).Value = 668
$ws.Cells.Item(30, 9).Value = 32
$ws.Cells.Item(31, 9).Value = 91
$ws.Cells.Item(33, 9).Value = 465
$ws.Cells.Item(36, 9).Value = 133
$ws.Cells.Item(37, 9).Value = 323
$ws.Cells.Item(38, 9).Value = 5
$ws.Cells.Item(42, 9).Value = 342
$ws.Cells.Item(44, 9).Value = 77
$ws.Cells.Item(50, 9).Value = 40
$ws.Cells.Item(52, 9).Value = 209
$ws.Cells.Item(55, 9).Value = 107
$ws.Cells.Item(57, 9).Value = 35
$ws.Cells.Item(63, 5).Value = 310
$ws.Cells.Item(65, 9).Value = 221
$ws.Cells.Item(67, 9).Value = 391
$ws.Cells.Item(68, 9).Value = 32
$ws.Cells.Item(70, 9).Value = 19
$ws.Cells.Item(71, 9).Value = 26
$ws.Cells.Item(72, 9).Value = 36
$ws.Cells.Item(76, 9).Value = 158
$ws.Cells.Item(78, 9).Value = 136
$ws.Cells.Item(83, 9).Value = 201
$ws.Cells.Item(85, 9).Value = 461
$ws.Cells.Item(90, 9).Value = 118
$ws.Cells.Item(92, 9).Value = 31
$ws.Cells.Item(94, 9).Value = 89
$ws.Cells.Item(95, 9).Value = 162
$ws.Cells.Item(99, 9).Value = 179
$ws.Cells.Item(101, 5).Value = 25967
$ws.Cells.Item(101, 9).Value = 9986

# South Shore (sheet 3)
$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 9).Value = 117
$ws.Cells.Item(6, 9).Value = 118
$ws.Cells.Item(7, 9).Value = 461

# Little Village (sheet 5)
$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(2, 9).Value = 52
$ws.Cells.Item(7, 9).Value = 209

# Austin (sheet 7)
$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 9).Value = 198
$ws.Cells.Item(6, 9).Value = 203
$ws.Cells.Item(7, 9).Value = 635

# Auburn Gresham (sheet 9)
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 9).Value = 116
$ws.Cells.Item(3, 9).Value = 101
$ws.Cells.Item(6, 9).Value = 85
$ws.Cells.Item(7, 9).Value = 335

# Bridgeport (sheet 12)
$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Cells.Item(6, 9).Value = 22
$ws.Cells.Item(7, 9).Value = 51

# Fuller Park (sheet 13)
$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Cells.Item(6, 9).Value = 9
$ws.Cells.Item(7, 9).Value = 32

# Grand Crossing (sheet 14)
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 9).Value = 104
$ws.Cells.Item(7, 9).Value = 323

# Woodlawn (sheet 15)
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(2, 9).Value = 44
$ws.Cells.Item(6, 9).Value = 52
$ws.Cells.Item(7, 9).Value = 179

# North Lawndale (sheet 16)
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(2, 9).Value = 90
$ws.Cells.Item(6, 9).Value = 139
$ws.Cells.Item(7, 9).Value = 391

# Gage Park (sheet 17)
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(2, 9).Value = 29
$ws.Cells.Item(7, 9).Value = 91

# New City (sheet 19)
$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(2, 9).Value = 71
$ws.Cells.Item(7, 9).Value = 221

# South Chicago (sheet 20)
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(2, 9).Value = 73
$ws.Cells.Item(3, 9).Value = 77
$ws.Cells.Item(7, 9).Value = 201

# West Pullman (sheet 21)
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(3, 9).Value = 63
$ws.Cells.Item(7, 9).Value = 162

# Garfield Park (sheet 22)
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 9).Value = 107
$ws.Cells.Item(3, 9).Value = 169
$ws.Cells.Item(6, 9).Value = 151
$ws.Cells.Item(7, 9).Value = 465

# Englewood (sheet 25)
$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 9).Value = 207
$ws.Cells.Item(3, 9).Value = 232
$ws.Cells.Item(5, 9).Value = 24
$ws.Cells.Item(6, 9).Value = 182
$ws.Cells.Item(7, 9).Value = 668

# Chatham (sheet 26)
$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(2, 9).Value = 113
$ws.Cells.Item(7, 9).Value = 274

# Irving Park (sheet 27)
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Cells.Item(6, 9).Value = 23
$ws.Cells.Item(7, 9).Value = 77

# River North (sheet 29)
$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(6, 9).Value = 64
$ws.Cells.Item(7, 9).Value = 158

# Humboldt Park (sheet 32)
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 9).Value = 91
$ws.Cells.Item(7, 9).Value = 342

# Rogers Park (sheet 35)
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(2, 9).Value = 26
$ws.Cells.Item(4, 9).Value = 21
$ws.Cells.Item(6, 9).Value = 55
$ws.Cells.Item(7, 9).Value = 136

# Lower West Side (sheet 36)
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(3, 9).Value = 29
$ws.Cells.Item(7, 9).Value = 107

# Douglas (sheet 39)
$ws = $wb.Worksheets.Item('Douglas')
$ws.Cells.Item(6, 9).Value = 24
$ws.Cells.Item(7, 9).Value = 87

# Chicago Lawn (sheet 44)
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(4, 9).Value = 14
$ws.Cells.Item(7, 9).Value = 252

# Grand Boulevard (sheet 47)
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(6, 9).Value = 43
$ws.Cells.Item(7, 9).Value = 133

# West Loop (sheet 51)
$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(6, 9).Value = 48
$ws.Cells.Item(7, 9).Value = 89

# Lincoln Square (sheet 56)
$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Cells.Item(6, 9).Value = 12
$ws.Cells.Item(7, 9).Value = 40

# West Elsdon (sheet 66)
$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Cells.Item(6, 9).Value = 14
$ws.Cells.Item(7, 9).Value = 31

# O'Hare (sheet 67)
$ws = $wb.Worksheets.Item('O''Hare')
$ws.Cells.Item(2, 9).Value = 5
$ws.Cells.Item(7, 9).Value = 19

# Edgewater (sheet 71)
$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(4, 9).Value = 12
$ws.Cells.Item(7, 9).Value = 90

# Washington Heights (sheet 74)
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(3, 9).Value = 22
$ws.Cells.Item(6, 9).Value = 43
$ws.Cells.Item(7, 9).Value = 118

# North Park (sheet 76)
$ws = $wb.Worksheets.Item('North Park')
$ws.Cells.Item(5, 9).Value = 1
$ws.Cells.Item(7, 9).Value = 32

# Mckinley Park (sheet 77)
$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Cells.Item(2, 9).Value = 13
$ws.Cells.Item(7, 9).Value = 35

# Oakland (sheet 81)
$ws = $wb.Worksheets.Item('Oakland')
$ws.Cells.Item(2, 9).Value = 8
$ws.Cells.Item(7, 9).Value = 26

# Old Town (sheet 82)
$ws = $wb.Worksheets.Item('Old Town')
$ws.Cells.Item(4, 9).Value = 4
$ws.Cells.Item(7, 9).Value = 36

# Grant Park (sheet 100)
$ws = $wb.Worksheets.Item('Grant Park')
$ws.Cells.Item(2, 9).Value = 2
$ws.Cells.Item(6, 9).Value = 5
